# Auto-generated Excel COM-interop script
# Adds 2022-08-15 daily crime counts to the "2022" (column I) totals
# across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('I2').Value = 4423
$ws.Range('I3').Value = 4626
$ws.Range('I4').Value = 1066
$ws.Range('I5').Value = 423
$ws.Range('I6').Value = 5043
$ws.Range('I7').Value = 15581

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('I6').Value = 64
$ws.Range('I7').Value = 179

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('I3').Value = 46
$ws.Range('I7').Value = 162

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('I2').Value = 31
$ws.Range('I7').Value = 84

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('I2').Value = 145
$ws.Range('I7').Value = 611

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('I6').Value = 52
$ws.Range('I7').Value = 149

$ws = $wb.Worksheets.Item('New City')
$ws.Range('I2').Value = 113
$ws.Range('I7').Value = 350

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('I2').Value = 126
$ws.Range('I4').Value = 58
$ws.Range('I8').Value = 946
$ws.Range('I10').Value = 109
$ws.Range('I11').Value = 238
$ws.Range('I12').Value = 35
$ws.Range('I14').Value = 84
$ws.Range('I15').Value = 178
$ws.Range('I16').Value = 41
$ws.Range('I19').Value = 439
$ws.Range('I20').Value = 373
$ws.Range('I29').Value = 995
$ws.Range('I31').Value = 149
$ws.Range('I33').Value = 723
$ws.Range('I34').Value = 73
$ws.Range('I36').Value = 218
$ws.Range('I48').Value = 215
$ws.Range('I50').Value = 72
$ws.Range('I51').Value = 169
$ws.Range('I60').Value = 77
$ws.Range('I63').Value = 60
$ws.Range('I65').Value = 350
$ws.Range('I67').Value = 611
$ws.Range('I73').Value = 133
$ws.Range('I76').Value = 232
$ws.Range('I78').Value = 224
$ws.Range('I79').Value = 431
$ws.Range('I85').Value = 696
$ws.Range('I86').Value = 92
$ws.Range('I89').Value = 179
$ws.Range('I90').Value = 192
$ws.Range('I91').Value = 185
$ws.Range('I94').Value = 145
$ws.Range('I95').Value = 256
$ws.Range('I96').Value = 162
$ws.Range('I97').Value = 120
$ws.Range('I98').Value = 101
$ws.Range('I101').Value = 15581

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('I3').Value = 99
$ws.Range('I6').Value = 45
$ws.Range('I7').Value = 256

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('I3').Value = 269
$ws.Range('I6').Value = 227
$ws.Range('I7').Value = 723

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('I3').Value = 70
$ws.Range('I4').Value = 25
$ws.Range('I6').Value = 167

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('I2').Value = 287
$ws.Range('I3').Value = 345
$ws.Range('I4').Value = 51
$ws.Range('I6').Value = 273
$ws.Range('I7').Value = 995

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('I6').Value = 120
$ws.Range('I7').Value = 439

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('I2').Value = 28
$ws.Range('I6').Value = 122
$ws.Range('I7').Value = 215

$ws = $wb.Worksheets.Item('River North')
$ws.Range('I2').Value = 49
$ws.Range('I7').Value = 232

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('I2').Value = 179
$ws.Range('I3').Value = 278
$ws.Range('I6').Value = 176
$ws.Range('I7').Value = 696

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('I6').Value = 46
$ws.Range('I7').Value = 109

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('I4').Value = 29
$ws.Range('I7').Value = 224

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('I3').Value = 68
$ws.Range('I7').Value = 185

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('I2').Value = 130
$ws.Range('I3').Value = 135
$ws.Range('I6').Value = 124
$ws.Range('I7').Value = 431

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('I4').Value = 28
$ws.Range('I7').Value = 373

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('I2').Value = 68
$ws.Range('I3').Value = 69
$ws.Range('I6').Value = 65
$ws.Range('I7').Value = 218

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('I2').Value = 30
$ws.Range('I7').Value = 73

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('I3').Value = 26
$ws.Range('I7').Value = 145

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('I2').Value = 54
$ws.Range('I6').Value = 64
$ws.Range('I7').Value = 178

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('I3').Value = 9
$ws.Range('I6').Value = 66
$ws.Range('I7').Value = 101

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('I4').Value = 14
$ws.Range('I5').Value = 1
$ws.Range('I7').Value = 72

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('I2').Value = 104
$ws.Range('I4').Value = 19
$ws.Range('I7').Value = 238

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('I3').Value = 43
$ws.Range('I7').Value = 133

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('I3').Value = 44
$ws.Range('I7').Value = 126

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('I6').Value = 72
$ws.Range('I7').Value = 120

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('I2').Value = 295
$ws.Range('I3').Value = 268
$ws.Range('I6').Value = 301
$ws.Range('I7').Value = 946

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('I4').Value = 44
$ws.Range('I7').Value = 92

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('I2').Value = 67
$ws.Range('I6').Value = 63
$ws.Range('I7').Value = 192

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('I6').Value = 68
$ws.Range('I7').Value = 169

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('I2').Value = 24
$ws.Range('I7').Value = 77

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range('I2').Value = 23
$ws.Range('I7').Value = 58

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range('I6').Value = 18
$ws.Range('I7').Value = 35

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('I6').Value = 26
$ws.Range('I7').Value = 41
